# Append rows 4-8 with bobina data (text-typed values, except numeric Sec on row 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new range to text format so numeric-looking strings
# (e.g. "12", "01") are preserved as text instead of being coerced to numbers.
$ws.Range("A4:K8").NumberFormat = "@"

# Row 4
$ws.Range("A4").Value = "12"
$ws.Range("B4").Value = "12"
$ws.Range("C4").Value = "21"
$ws.Range("D4").Value = "21"
$ws.Range("E4").Value = "122"
$ws.Range("F4").NumberFormat = "General"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "2122"
$ws.Range("H4").Value = "2025-02-06 23:18"
$ws.Range("I4").Value = "A"
$ws.Range("J4").Value = "01"
$ws.Range("K4").Value = "ONDA LINER"

# Row 5
$ws.Range("A5").Value = "12"
$ws.Range("B5").Value = "12"
$ws.Range("C5").Value = "21"
$ws.Range("D5").Value = "21"
$ws.Range("E5").Value = "122"
$ws.Range("F5").Value = "2"
$ws.Range("G5").Value = "2122"
$ws.Range("H5").Value = "2025-02-06 23:18"
$ws.Range("I5").Value = "A"
$ws.Range("J5").Value = "01"
$ws.Range("K5").Value = "ONDA LINER"

# Row 6
$ws.Range("A6").Value = "12"
$ws.Range("B6").Value = "12"
$ws.Range("C6").Value = "21"
$ws.Range("D6").Value = "21"
$ws.Range("E6").Value = "122"
$ws.Range("F6").Value = "2"
$ws.Range("G6").Value = "2122"
$ws.Range("H6").Value = "2025-02-06 23:18"
$ws.Range("I6").Value = "A"
$ws.Range("J6").Value = "01"
$ws.Range("K6").Value = "ONDA LINER"

# Row 7
$ws.Range("A7").Value = "12"
$ws.Range("B7").Value = "12"
$ws.Range("C7").Value = "21"
$ws.Range("D7").Value = "21"
$ws.Range("E7").Value = "122"
$ws.Range("F7").Value = "3"
$ws.Range("G7").Value = "2122"
$ws.Range("H7").Value = "2025-02-06 23:18"
$ws.Range("I7").Value = "A"
$ws.Range("J7").Value = "01"
$ws.Range("K7").Value = "ONDA LINER"

# Row 8
$ws.Range("A8").Value = "12"
$ws.Range("B8").Value = "12"
$ws.Range("C8").Value = "21"
$ws.Range("D8").Value = "21"
$ws.Range("E8").Value = "122"
$ws.Range("F8").Value = "3"
$ws.Range("G8").Value = "2122"
$ws.Range("H8").Value = "2025-02-06 23:18"
$ws.Range("I8").Value = "A"
$ws.Range("J8").Value = "01"
$ws.Range("K8").Value = "ONDA LINER"

